$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-04-05 Saturday" "2025-04-06 Sunday"

Replace-Text "221÷7=" "255÷7="
Replace-Text "552÷5=" "287÷8="
Replace-Text "285÷5=" "487÷7="
Replace-Text "993÷9=" "507÷5="
Replace-Text "798÷6=" "445÷3="

Replace-Text "843÷2=" "475÷5="
Replace-Text "355÷9=" "167÷6="
Replace-Text "561÷3=" "427÷2="
Replace-Text "647÷3=" "978÷2="
Replace-Text "890÷5=" "517÷7="

Replace-Text "292÷3=" "283÷8="
Replace-Text "770÷3=" "158÷5="
Replace-Text "283÷2=" "322÷3="
Replace-Text "231÷9=" "763÷2="
Replace-Text "837÷5=" "768÷7="

Replace-Text "290÷2=" "976÷4="
Replace-Text "807÷2=" "914÷6="
Replace-Text "773÷3=" "930÷9="
Replace-Text "690÷7=" "198÷3="
Replace-Text "274÷4=" "702÷5="

Replace-Text "476÷3=" "430÷8="
Replace-Text "867÷2=" "702÷8="
Replace-Text "268÷9=" "418÷7="
Replace-Text "885÷4=" "869÷3="
Replace-Text "684÷7=" "663÷5="
